# Append the latest batch of participants (sub_045..sub_047) to the QC
# debrief/errors check sheet, each starting out as "not an error"
# (FALSE in the boolean "reason" flag column), mirroring the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipants = @("sub_045", "sub_046", "sub_047")

$startRow = 46
for ($i = 0; $i -lt $newParticipants.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newParticipants[$i]
    $ws.Cells.Item($row, 2).Value = $false
}

# Leave the view looking at the newly added rows, same as the author did.
$ws.Range("A41:B48").Select()
